# Applies: Dataframe operations -> new "Salary" column + 3 new employees
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new employee names first so the shared-string table order matches
# (Hajra, Nashra, Zainab, then Salary).
$ws.Range("A7").Value = "Hajra"
$ws.Range("A8").Value = "Nashra"
$ws.Range("A9").Value = "Zainab"

# New column D header: Salary
$ws.Range("D1").Value = "Salary"

# Salary values for existing rows
$ws.Range("D2").Value = 1200000
$ws.Range("D4").Value = 1000000
$ws.Range("D5").Value = 600000
$ws.Range("D6").Value = 500000

# Remove Zishan's Age value (B4) - clear just that cell
$ws.Range("B4").ClearContents()

# Remaining new-row data
$ws.Range("B7").Value = 24
$ws.Range("D7").Value = 450000

$ws.Range("B8").Value = 2

$ws.Range("B10").Value = 36
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 3500000

$ws.Range("B11").Value = 38

# Update selection to match target state
$ws.Range("H7:I10").Select()
